$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Shorten three existing "Loc" keywords (drop the trailing word) ---
#     A31: "living room"  -> "livingroom"
#     A33: "dining room"  -> "dining"
#     A36: "laundry room" -> "laundry"
$ws.Range("A31").Value = "livingroom"
$ws.Range("A33").Value = "dining"
$ws.Range("A36").Value = "laundry"

# --- 2. Add three new keyword rows at the bottom of the "Loc" block ---
#     row150: room
#     row151: living
#     row152: living room
# Each follows the same layout as the other Loc rows: column B (Loc) = 1,
# columns C:H = 0.
$styleAF = $ws.Range("A2").Style
$styleGH = $ws.Range("G2").Style

$ws.Range("A150").Value = "room"
$ws.Range("B150").Value = 1
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 0
$ws.Range("A150:F150").Style = $styleAF
$ws.Range("G150:H150").Style = $styleGH

$ws.Range("A151").Value = "living"
$ws.Range("B151").Value = 1
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0
$ws.Range("A151:F151").Style = $styleAF
$ws.Range("G151:H151").Style = $styleGH

$ws.Range("A152").Value = "living room"
$ws.Range("B152").Value = 1
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 0
$ws.Range("E152").Value = 0
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 0
$ws.Range("A152:F152").Style = $styleAF
$ws.Range("G152:H152").Style = $styleGH

# --- 3. Restore the saved selection/active cell ---
$ws.Range("M132").Select()
